$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 9).Value = 'sd'
$ws.Cells.Item(7, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(8, 9).Value = 'aa'
$ws.Cells.Item(8, 10).Value = 'Agree/Accept'
$ws.Cells.Item(9, 9).Value = 'sd'
$ws.Cells.Item(9, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(12, 9).Value = 'ba'
$ws.Cells.Item(12, 10).Value = 'Appreciation'
$ws.Cells.Item(14, 9).Value = 'b'
$ws.Cells.Item(14, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(16, 9).Value = 'sd'
$ws.Cells.Item(16, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(18, 9).Value = 'ba'
$ws.Cells.Item(18, 10).Value = 'Appreciation'
$ws.Cells.Item(19, 9).Value = 'b'
$ws.Cells.Item(19, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(25, 9).Value = 'sv'
$ws.Cells.Item(25, 10).Value = 'Statement-opinion'
$ws.Cells.Item(29, 9).Value = 'sd'
$ws.Cells.Item(29, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(33, 9).Value = 'aa'
$ws.Cells.Item(33, 10).Value = 'Agree/Accept'
$ws.Cells.Item(34, 9).Value = 'sv'
$ws.Cells.Item(34, 10).Value = 'Statement-opinion'
$ws.Cells.Item(38, 9).Value = 'sd'
$ws.Cells.Item(38, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(49, 9).Value = 'sv'
$ws.Cells.Item(49, 10).Value = 'Statement-opinion'
$ws.Cells.Item(65, 9).Value = 'ba'
$ws.Cells.Item(65, 10).Value = 'Appreciation'
$ws.Cells.Item(67, 9).Value = 'sd'
$ws.Cells.Item(67, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(72, 9).Value = 'aa'
$ws.Cells.Item(72, 10).Value = 'Agree/Accept'
$ws.Cells.Item(74, 9).Value = 'aa'
$ws.Cells.Item(74, 10).Value = 'Agree/Accept'
$ws.Cells.Item(75, 9).Value = 'aa'
$ws.Cells.Item(75, 10).Value = 'Agree/Accept'
$ws.Cells.Item(77, 9).Value = 'ba'
$ws.Cells.Item(77, 10).Value = 'Appreciation'
$ws.Cells.Item(86, 9).Value = 'ba'
$ws.Cells.Item(86, 10).Value = 'Appreciation'
$ws.Cells.Item(88, 9).Value = 'aa'
$ws.Cells.Item(88, 10).Value = 'Agree/Accept'
$ws.Cells.Item(94, 9).Value = 'aa'
$ws.Cells.Item(94, 10).Value = 'Agree/Accept'
$ws.Cells.Item(115, 9).Value = 'sv'
$ws.Cells.Item(115, 10).Value = 'Statement-opinion'
$ws.Cells.Item(121, 9).Value = '%'
$ws.Cells.Item(121, 10).Value = 'Uninterpretable'
$ws.Cells.Item(125, 9).Value = 'sv'
$ws.Cells.Item(125, 10).Value = 'Statement-opinion'
$ws.Cells.Item(133, 9).Value = 'sd'
$ws.Cells.Item(133, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(136, 9).Value = 'ba'
$ws.Cells.Item(136, 10).Value = 'Appreciation'
$ws.Cells.Item(137, 9).Value = 'ba'
$ws.Cells.Item(137, 10).Value = 'Appreciation'
$ws.Cells.Item(138, 9).Value = 'ba'
$ws.Cells.Item(138, 10).Value = 'Appreciation'
$ws.Cells.Item(139, 9).Value = 'ba'
$ws.Cells.Item(139, 10).Value = 'Appreciation'
$ws.Cells.Item(141, 9).Value = 'b'
$ws.Cells.Item(141, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(150, 9).Value = 'sv'
$ws.Cells.Item(150, 10).Value = 'Statement-opinion'
$ws.Cells.Item(151, 9).Value = 'sv'
$ws.Cells.Item(151, 10).Value = 'Statement-opinion'
$ws.Cells.Item(157, 9).Value = 'aa'
$ws.Cells.Item(157, 10).Value = 'Agree/Accept'
$ws.Cells.Item(167, 9).Value = 'sd'
$ws.Cells.Item(167, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(170, 9).Value = 'sv'
$ws.Cells.Item(170, 10).Value = 'Statement-opinion'
$ws.Cells.Item(172, 9).Value = '%'
$ws.Cells.Item(172, 10).Value = 'Uninterpretable'
$ws.Cells.Item(179, 9).Value = 'sd'
$ws.Cells.Item(179, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(183, 9).Value = 'sv'
$ws.Cells.Item(183, 10).Value = 'Statement-opinion'
$ws.Cells.Item(191, 9).Value = 'qy'
$ws.Cells.Item(191, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(193, 9).Value = 'sv'
$ws.Cells.Item(193, 10).Value = 'Statement-opinion'
$ws.Cells.Item(199, 9).Value = 'sd'
$ws.Cells.Item(199, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(201, 9).Value = 'ba'
$ws.Cells.Item(201, 10).Value = 'Appreciation'
$ws.Cells.Item(217, 9).Value = 'sd'
$ws.Cells.Item(217, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(226, 9).Value = 'sd'
$ws.Cells.Item(226, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(227, 9).Value = 'b'
$ws.Cells.Item(227, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(247, 9).Value = 'ba'
$ws.Cells.Item(247, 10).Value = 'Appreciation'
$ws.Cells.Item(254, 9).Value = 'sd'
$ws.Cells.Item(254, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(262, 9).Value = 'sv'
$ws.Cells.Item(262, 10).Value = 'Statement-opinion'
$ws.Cells.Item(278, 9).Value = 'sv'
$ws.Cells.Item(278, 10).Value = 'Statement-opinion'
$ws.Cells.Item(280, 9).Value = 'sd'
$ws.Cells.Item(280, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(283, 9).Value = '%'
$ws.Cells.Item(283, 10).Value = 'Uninterpretable'
$ws.Cells.Item(295, 9).Value = 'sd'
$ws.Cells.Item(295, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(308, 9).Value = 'sd'
$ws.Cells.Item(308, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(329, 9).Value = 'b'
$ws.Cells.Item(329, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(347, 9).Value = 'sd'
$ws.Cells.Item(347, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(352, 9).Value = 'sd'
$ws.Cells.Item(352, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(364, 9).Value = 'sv'
$ws.Cells.Item(364, 10).Value = 'Statement-opinion'
$ws.Cells.Item(371, 9).Value = 'aa'
$ws.Cells.Item(371, 10).Value = 'Agree/Accept'
$ws.Cells.Item(378, 9).Value = 'ba'
$ws.Cells.Item(378, 10).Value = 'Appreciation'
$ws.Cells.Item(379, 9).Value = 'aa'
$ws.Cells.Item(379, 10).Value = 'Agree/Accept'
$ws.Cells.Item(391, 9).Value = 'sd'
$ws.Cells.Item(391, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(404, 9).Value = 'sv'
$ws.Cells.Item(404, 10).Value = 'Statement-opinion'
$ws.Cells.Item(405, 9).Value = 'ba'
$ws.Cells.Item(405, 10).Value = 'Appreciation'
$ws.Cells.Item(414, 9).Value = 'b'
$ws.Cells.Item(414, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(419, 9).Value = 'ba'
$ws.Cells.Item(419, 10).Value = 'Appreciation'
$ws.Cells.Item(427, 9).Value = 'sv'
$ws.Cells.Item(427, 10).Value = 'Statement-opinion'
$ws.Cells.Item(430, 9).Value = 'ba'
$ws.Cells.Item(430, 10).Value = 'Appreciation'
$ws.Cells.Item(464, 9).Value = 'sv'
$ws.Cells.Item(464, 10).Value = 'Statement-opinion'
$ws.Cells.Item(472, 9).Value = 'aa'
$ws.Cells.Item(472, 10).Value = 'Agree/Accept'
$ws.Cells.Item(481, 9).Value = 'sv'
$ws.Cells.Item(481, 10).Value = 'Statement-opinion'
$ws.Cells.Item(487, 9).Value = 'sd'
$ws.Cells.Item(487, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(489, 9).Value = 'sd'
$ws.Cells.Item(489, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(492, 9).Value = 'ba'
$ws.Cells.Item(492, 10).Value = 'Appreciation'
$ws.Cells.Item(500, 9).Value = 'sd'
$ws.Cells.Item(500, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(501, 9).Value = 'sd'
$ws.Cells.Item(501, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(502, 9).Value = 'ba'
$ws.Cells.Item(502, 10).Value = 'Appreciation'
$ws.Cells.Item(522, 9).Value = 'sd'
$ws.Cells.Item(522, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(524, 9).Value = 'sd'
$ws.Cells.Item(524, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(525, 9).Value = 'sv'
$ws.Cells.Item(525, 10).Value = 'Statement-opinion'
